# Update column C (row 2 through row 45) from 45221 (2023-10-22) to 45224 (2023-10-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
